$d = $word.ActiveDocument

# Update the date/day heading (text is unique in the document, Find is safe here)
$d.Content.Find.Execute("2024-03-20 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-03-21 Thursday", 2)

# Update the division problems in the table, cell by cell (row, col), using
# direct Range.Text assignment. Several of the problem strings repeat across
# different cells, so Find/Replace (even scoped to a cell Range) is not safe
# here -- it matches across the whole document. Direct assignment to the
# cell's own Range correctly targets only that cell, and preserves the
# run/paragraph formatting already in place.
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $table.Cell($row, $col).Range.Text = $newText
}

# Row 1 (table row 1)
Set-CellText $t 1 1 "55÷3="
Set-CellText $t 1 2 "89÷7="
Set-CellText $t 1 3 "83÷9="
Set-CellText $t 1 4 "46÷6="
Set-CellText $t 1 5 "61÷8="

# Row 2 (table row 5)
Set-CellText $t 5 1 "84÷4="
Set-CellText $t 5 2 "13÷9="
Set-CellText $t 5 3 "44÷7="
Set-CellText $t 5 4 "51÷2="
Set-CellText $t 5 5 "82÷7="

# Row 3 (table row 9)
Set-CellText $t 9 1 "59÷9="
Set-CellText $t 9 2 "80÷2="
Set-CellText $t 9 3 "41÷5="
Set-CellText $t 9 4 "77÷3="
Set-CellText $t 9 5 "34÷5="

# Row 4 (table row 13)
Set-CellText $t 13 1 "69÷8="
Set-CellText $t 13 2 "44÷8="
Set-CellText $t 13 3 "86÷6="
Set-CellText $t 13 4 "55÷5="
Set-CellText $t 13 5 "60÷2="

# Row 5 (table row 17)
Set-CellText $t 17 1 "46÷6="
Set-CellText $t 17 2 "27÷5="
Set-CellText $t 17 3 "44÷8="
Set-CellText $t 17 4 "74÷5="
Set-CellText $t 17 5 "89÷4="
